$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy existing formatting (styles) from row 21 down into the new/changed
# date and time cells before setting values, so no new style entries are
# created in styles.xml (mirrors the original author's formatting reuse).
$ws.Range("A21").Copy()
$ws.Range("A23:A24").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B21:C21").Copy()
$ws.Range("B22:C24").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Row 22: fill in B22:E22 (A22 and F22 already existed)
$ws.Range("B22").Value = 0.33333333333333331
$ws.Range("C22").Value = 0.46875
$ws.Range("D22").Value = "Juliano, Tommy, Constantin, Jot"
$ws.Range("E22").Value = "Suche nach Lösung für verschiedene Displaygrößen"

# Row 23: fill in A23:E23 (F23 already existed)
$ws.Range("A23").Value = 42110
$ws.Range("B23").Value = 0.47916666666666669
$ws.Range("C23").Value = 0.54166666666666663
$ws.Range("D23").Value = "Constantin, Juliano"
$ws.Range("E23").Value = "Einbau einer Buttonanimation beim Klick"

# Row 24: new row entirely
$ws.Range("A24").Value = 42112
$ws.Range("B24").Value = 0.58333333333333337
$ws.Range("C24").Value = 0.75
$ws.Range("D24").Value = "Juliano"
$ws.Range("E24").Value = "Anlegen der drawable-Ordner für verschiedene Displaygrößen und rendern der Hintergrundbilder"

# Update selection to mimic final cursor position from the diff
$ws.Range("E24").Select()
